$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.86"
$ws.Range("E2").Value = "'1.86%"

$ws.Range("D3").Value = "'27.21"
$ws.Range("E3").Value = "'1.86%"

$ws.Range("D4").Value = "'4.737"
$ws.Range("E4").Value = "'8.12%"

$ws.Range("D5").Value = "'0.06077"
$ws.Range("E5").Value = "'3.51%"

$ws.Range("D6").Value = "'6.666"
$ws.Range("E6").Value = "'0.85%"

$ws.Range("D7").Value = "'0.8474"
$ws.Range("E7").Value = "'-0.52%"

$ws.Range("D8").Value = "'0.9173"
$ws.Range("E8").Value = "'-1.19%"

$ws.Range("D10").Value = "'0.04933"
$ws.Range("E10").Value = "'7.65%"

$ws.Range("D11").Value = "'0.07084"
$ws.Range("E11").Value = "'0.70%"

$ws.Range("D12").Value = "'0.03133"
$ws.Range("E12").Value = "'2.13%"

$ws.Range("D13").Value = "'0.09076"
$ws.Range("E13").Value = "'-0.29%"

$ws.Range("D14").Value = "'0.001531"
$ws.Range("E14").Value = "'-0.86%"

$ws.Range("D15").Value = "'0.0006096"
$ws.Range("E15").Value = "'-94.09%"

$ws.Range("D16").Value = "'0.006128"
$ws.Range("E16").Value = "'1.52%"

$ws.Range("D18").Value = "'3.154"
$ws.Range("E18").Value = "'-0.68%"

$ws.Range("E20").Value = "'2.54%"

$ws.Range("D21").Value = "'0.1297"
$ws.Range("E21").Value = "'0.86%"

$ws.Range("D22").Value = "'4.095"
$ws.Range("E22").Value = "'4.77%"

$ws.Range("D23").Value = "'0.04255"
$ws.Range("E23").Value = "'-0.07%"

$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-0.06%"

$ws.Range("E25").Value = "'-8.84%"

$ws.Range("E26").Value = "'-0.02%"

$ws.Range("E27").Value = "'3.08%"

$ws.Range("D40").Value = "'0.03875"
$ws.Range("E40").Value = "'1.88%"

$ws.Range("E41").Value = "'1.36%"

$ws.Range("D42").Value = "'0.004130"
$ws.Range("E42").Value = "'-33.84%"

$ws.Range("D43").Value = "'0.01635"
$ws.Range("E43").Value = "'18.06%"

$ws.Range("E44").Value = "'0.34%"

$ws.Range("D45").Value = "'0.00005330"
$ws.Range("E45").Value = "'-0.92%"

$ws.Range("E46").Value = "'-0.02%"

$ws.Range("E47").Value = "'1.21%"

$ws.Range("D48").Value = "'0.1353"
$ws.Range("E48").Value = "'-46.24%"

$ws.Range("E49").Value = "'-0.02%"

$ws.Range("E50").Value = "'-0.02%"
